$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "GT2 Pulley 16T 5mm" line item (row 19) entirely, shifting
# everything below up by one row.
$ws.Rows("19").Delete()

# Remove the "500mm Pulley Belt" line item (originally row 26, now row 25
# after the previous delete) entirely.
$ws.Rows("25").Delete()

# Antenna and Xbee turned out to be free / donated -- zero out their unit
# cost (and thus their computed total).
$ws.Range("E21").Value = 0
$ws.Range("E22").Value = 0

# Add the new "12x12 ABS Sheet" purchase as the new last row of the table.
$ws.Range("A29").Value = 45198
$ws.Range("C29").Value = "12x12 ABS Sheet"
$ws.Range("D29").Value = 3
$ws.Range("E29").Value = 4.26
$ws.Range("F29").Formula = "=D29*E29"
$ws.Range("G29").Value = "Alex Treseder"
